$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Heat Storage")

# --- Shared-string slot trick -------------------------------------------------
# The sole current user of the "OPEX fix [€/MWha]" shared string is I1.
# Renaming that cell's text first causes the engine to edit that shared
# string *in place* (since it becomes the only / soon-to-be-freed user),
# which is what we want for the brand new "Hourly Stoarge Losses [%]"
# header text (it reuses that slot). We fix the real OPEX label afterwards,
# once the insert has shifted this cell out of the way, so it gets a new
# shared-string slot appended at the end.
$ws.Range("I1").Value = "Hourly Stoarge Losses [%]"

# --- Insert a new column for "Hourly Stoarge Losses [%]" ---------------------
$ws.Columns.Item(4).Insert()

# The old I1 cell (now shifted to J1) currently holds "Hourly Stoarge
# Losses [%]" - move that text to the new D1 cell, and restore the correct
# (typo-fixed) OPEX header text at J1.
$ws.Range("D1").Value = "Hourly Stoarge Losses [%]"
$ws.Range("J1").Value = "OPEX fix [€/MWh]"

# --- Header formatting ---------------------------------------------------
# Make every header cell B1:K1 look like the existing bold/filled/bordered
# header style (same style already used by H1/I1 before our edits).
$ws.Range("I1").Copy()
$ws.Range("B1:K1").PasteSpecial(-4122)   # xlPasteFormats

# E1 ("maximum unloading power [MW]") additionally gets word-wrap since its
# column is narrower now.
$ws.Range("E1").WrapText = $true

$ws.Rows.Item(1).RowHeight = 30

# --- Row 2 values ----------------------------------------------------------
# New column D2 (Hourly Stoarge Losses) is left blank - no data supplied.
# Investment cost (now I2) changes value from 60 to 3000.
$ws.Range("I2").Value = 3000

# Clear the leftover border-only style from the cells that used to carry it
# (old H2/I2) now that the data has shifted - row 2 data cells are unstyled.
$ws.Range("I2:J2").ClearFormats()
$ws.Range("I2").Value = 3000
$ws.Range("J2").Value = 10000

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 2.14
$ws.Columns.Item(4).ColumnWidth = 24.14
$ws.Columns.Item(5).ColumnWidth = 25.86
$ws.Columns.Item(10).ColumnWidth = 17.57

# --- Selection -----------------------------------------------------------
$ws.Range("J1").Select()
